$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 holds the most recent day's spot price data. This automated update
# refreshes the date and all hourly/slot values for the new day.

$ws.Range("A2").Value = 46026

$ws.Range("B2").Value = 87.73999999999999
$ws.Range("C2").Value = 65.93000000000001
$ws.Range("D2").Value = 60.11
$ws.Range("E2").Value = 57.57
$ws.Range("F2").Value = 53.38
$ws.Range("G2").Value = 50.43
$ws.Range("H2").Value = 54.39
$ws.Range("I2").Value = 56.76
$ws.Range("J2").Value = 57.79
$ws.Range("K2").Value = 55.96
$ws.Range("L2").Value = 55.55
$ws.Range("M2").Value = 57.98
$ws.Range("N2").Value = 58.11
$ws.Range("O2").Value = 60.09
$ws.Range("P2").Value = 61.24
$ws.Range("Q2").Value = 61.01
$ws.Range("R2").Value = 62.68
$ws.Range("S2").Value = 68.09
$ws.Range("T2").Value = 80.84
$ws.Range("U2").Value = 84.56999999999999
$ws.Range("V2").Value = 88.18000000000001
$ws.Range("W2").Value = 97.68000000000001
$ws.Range("X2").Value = 87.95999999999999
$ws.Range("Y2").Value = 78.76000000000001
$ws.Range("Z2").Value = 66.78

$ws.Range("AA2").Value = "20h-24h"
$ws.Range("AB2").Value = 88.15000000000001
$ws.Range("AC2").Value = "20h-22h"
$ws.Range("AD2").Value = 92.93000000000001
$ws.Range("AE2").Value = "22h-24h"
$ws.Range("AF2").Value = 83.36
$ws.Range("AG2").Value = "1h-16h"
